# Natmi following Dr Hou advice
#
# The NATMI ligand-receptor (Lama2-Itga7) edge table is recomputed to add a
# new "M2" (macrophage) cluster into both the sending-cluster and
# target-cluster dimensions. The table is a full cross-product of
# Sending cluster x Target cluster (previously ECs/FAPs/sCs = 3 senders x 4
# targets = 12 data rows; now ECs/FAPs/M2/sCs = 4 senders x 4 targets = 16
# data rows), so the sheet grows from A1:T13 (header + 12 data rows) to
# A1:T17 (header + 16 data rows) and all of the previously-computed
# statistics shift because the underlying per-cluster expression pools
# changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lama2"
$ws.Range("C2").Value = "Itga7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.608931666666666
$ws.Range("H2").Value = 13.826795
$ws.Range("I2").Value = 0.02269509467890621
$ws.Range("J2").Value = 0.02269509467890622
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.636552333333334
$ws.Range("N2").Value = 16.909657
$ws.Range("O2").Value = 0.09690140221006956
$ws.Range("P2").Value = 0.09690140221006956
$ws.Range("Q2").Value = 25.97848453992389
$ws.Range("R2").Value = 233.806360859315
$ws.Range("S2").Value = 0.0021991864976763
$ws.Range("T2").Value = 0.002199186497676301
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lama2"
$ws.Range("C3").Value = "Itga7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.608931666666666
$ws.Range("H3").Value = 13.826795
$ws.Range("I3").Value = 0.02269509467890621
$ws.Range("J3").Value = 0.02269509467890622
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.050135999999999
$ws.Range("N3").Value = 12.150408
$ws.Range("O3").Value = 0.06962835335006774
$ws.Range("P3").Value = 0.06962835335006774
$ws.Range("Q3").Value = 18.66680006470666
$ws.Range("R3").Value = 168.00120058236
$ws.Range("S3").Value = 0.001580222071616124
$ws.Range("T3").Value = 0.001580222071616124
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lama2"
$ws.Range("C4").Value = "Itga7"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.608931666666666
$ws.Range("H4").Value = 13.826795
$ws.Range("I4").Value = 0.02269509467890621
$ws.Range("J4").Value = 0.02269509467890622
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6741636666666667
$ws.Range("N4").Value = 2.022491
$ws.Range("O4").Value = 0.01158995796645939
$ws.Range("P4").Value = 0.01158995796645939
$ws.Range("Q4").Value = 3.107174271816111
$ws.Range("R4").Value = 27.964568446345
$ws.Range("S4").Value = 0.0002630351933733392
$ws.Range("T4").Value = 0.0002630351933733392
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lama2"
$ws.Range("C5").Value = "Itga7"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.608931666666666
$ws.Range("H5").Value = 13.826795
$ws.Range("I5").Value = 0.02269509467890621
$ws.Range("J5").Value = 0.02269509467890622
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 47.807061
$ws.Range("N5").Value = 143.421183
$ws.Range("O5").Value = 0.8218802864734033
$ws.Range("P5").Value = 0.8218802864734033
$ws.Range("Q5").Value = 220.339477333165
$ws.Range("R5").Value = 1983.055295998485
$ws.Range("S5").Value = 0.01865265091624045
$ws.Range("T5").Value = 0.01865265091624045
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama2"
$ws.Range("C6").Value = "Itga7"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 168.218394
$ws.Range("H6").Value = 504.655182
$ws.Range("I6").Value = 0.8283334739316415
$ws.Range("J6").Value = 0.8283334739316416
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.636552333333334
$ws.Range("N6").Value = 16.909657
$ws.Range("O6").Value = 0.09690140221006956
$ws.Range("P6").Value = 0.09690140221006956
$ws.Range("Q6").Value = 948.1717812102862
$ws.Range("R6").Value = 8533.546030892576
$ws.Range("S6").Value = 0.08026667512151416
$ws.Range("T6").Value = 0.08026667512151417
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lama2"
$ws.Range("C7").Value = "Itga7"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 168.218394
$ws.Range("H7").Value = 504.655182
$ws.Range("I7").Value = 0.8283334739316415
$ws.Range("J7").Value = 0.8283334739316416
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.050135999999999
$ws.Range("N7").Value = 12.150408
$ws.Range("O7").Value = 0.06962835335006774
$ws.Range("P7").Value = 0.06962835335006774
$ws.Range("Q7").Value = 681.307373401584
$ws.Range("R7").Value = 6131.766360614256
$ws.Range("S7").Value = 0.05767549581460146
$ws.Range("T7").Value = 0.05767549581460147
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lama2"
$ws.Range("C8").Value = "Itga7"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 168.218394
$ws.Range("H8").Value = 504.655182
$ws.Range("I8").Value = 0.8283334739316415
$ws.Range("J8").Value = 0.8283334739316416
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6741636666666667
$ws.Range("N8").Value = 2.022491
$ws.Range("O8").Value = 0.01158995796645939
$ws.Range("P8").Value = 0.01158995796645939
$ws.Range("Q8").Value = 113.406729299818
$ws.Range("R8").Value = 1020.660563698362
$ws.Range("S8").Value = 0.009600350145079008
$ws.Range("T8").Value = 0.009600350145079008
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lama2"
$ws.Range("C9").Value = "Itga7"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 168.218394
$ws.Range("H9").Value = 504.655182
$ws.Range("I9").Value = 0.8283334739316415
$ws.Range("J9").Value = 0.8283334739316416
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 47.807061
$ws.Range("N9").Value = 143.421183
$ws.Range("O9").Value = 0.8218802864734033
$ws.Range("P9").Value = 0.8218802864734033
$ws.Range("Q9").Value = 8042.027023280036
$ws.Range("R9").Value = 72378.24320952031
$ws.Range("S9").Value = 0.6807909528504468
$ws.Range("T9").Value = 0.6807909528504469
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Lama2"
$ws.Range("C10").Value = "Itga7"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1627236666666667
$ws.Range("H10").Value = 0.488171
$ws.Range("I10").Value = 0.00080127658394417
$ws.Range("J10").Value = 0.00080127658394417
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.636552333333334
$ws.Range("N10").Value = 16.909657
$ws.Range("O10").Value = 0.09690140221006956
$ws.Range("P10").Value = 0.09690140221006956
$ws.Range("Q10").Value = 0.9172004630385558
$ws.Range("R10").Value = 8.254804167347002
$ws.Range("S10").Value = 0.00007764482454228458
$ws.Range("T10").Value = 0.00007764482454228458
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Lama2"
$ws.Range("C11").Value = "Itga7"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1627236666666667
$ws.Range("H11").Value = 0.488171
$ws.Range("I11").Value = 0.00080127658394417
$ws.Range("J11").Value = 0.00080127658394417
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.050135999999999
$ws.Range("N11").Value = 12.150408
$ws.Range("O11").Value = 0.06962835335006774
$ws.Range("P11").Value = 0.06962835335006774
$ws.Range("Q11").Value = 0.6590529804186667
$ws.Range("R11").Value = 5.931476823768
$ws.Range("S11").Value = 0.00005579156911799989
$ws.Range("T11").Value = 0.00005579156911799989
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Lama2"
$ws.Range("C12").Value = "Itga7"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1627236666666667
$ws.Range("H12").Value = 0.488171
$ws.Range("I12").Value = 0.00080127658394417
$ws.Range("J12").Value = 0.00080127658394417
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6741636666666667
$ws.Range("N12").Value = 2.022491
$ws.Range("O12").Value = 0.01158995796645939
$ws.Range("P12").Value = 0.01158995796645939
$ws.Range("Q12").Value = 0.1097023837734445
$ws.Range("R12").Value = 0.987321453961
$ws.Range("S12").Value = 0.0000092867619274211
$ws.Range("T12").Value = 0.000009286761927421098
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Lama2"
$ws.Range("C13").Value = "Itga7"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1627236666666667
$ws.Range("H13").Value = 0.488171
$ws.Range("I13").Value = 0.00080127658394417
$ws.Range("J13").Value = 0.00080127658394417
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 47.807061
$ws.Range("N13").Value = 143.421183
$ws.Range("O13").Value = 0.8218802864734033
$ws.Range("P13").Value = 0.8218802864734033
$ws.Range("Q13").Value = 7.779340258477002
$ws.Range("R13").Value = 70.01406232629301
$ws.Range("S13").Value = 0.0006585534283564644
$ws.Range("T13").Value = 0.0006585534283564644
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Lama2"
$ws.Range("C14").Value = "Itga7"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 30.09047233333333
$ws.Range("H14").Value = 90.271417
$ws.Range("I14").Value = 0.148170154805508
$ws.Range("J14").Value = 0.148170154805508
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.636552333333334
$ws.Range("N14").Value = 16.909657
$ws.Range("O14").Value = 0.09690140221006956
$ws.Range("P14").Value = 0.09690140221006956
$ws.Range("Q14").Value = 169.6065220415522
$ws.Range("R14").Value = 1526.458698373969
$ws.Range("S14").Value = 0.01435789576633681
$ws.Range("T14").Value = 0.01435789576633681
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Lama2"
$ws.Range("C15").Value = "Itga7"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 30.09047233333333
$ws.Range("H15").Value = 90.271417
$ws.Range("I15").Value = 0.148170154805508
$ws.Range("J15").Value = 0.148170154805508
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.050135999999999
$ws.Range("N15").Value = 12.150408
$ws.Range("O15").Value = 0.06962835335006774
$ws.Range("P15").Value = 0.06962835335006774
$ws.Range("Q15").Value = 121.8705052542373
$ws.Range("R15").Value = 1096.834547288136
$ws.Range("S15").Value = 0.01031684389473215
$ws.Range("T15").Value = 0.01031684389473215
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Lama2"
$ws.Range("C16").Value = "Itga7"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 30.09047233333333
$ws.Range("H16").Value = 90.271417
$ws.Range("I16").Value = 0.148170154805508
$ws.Range("J16").Value = 0.148170154805508
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.6741636666666667
$ws.Range("N16").Value = 2.022491
$ws.Range("O16").Value = 0.01158995796645939
$ws.Range("P16").Value = 0.01158995796645939
$ws.Range("Q16").Value = 20.28590315997189
$ws.Range("R16").Value = 182.573128439747
$ws.Range("S16").Value = 0.001717285866079619
$ws.Range("T16").Value = 0.001717285866079619
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Lama2"
$ws.Range("C17").Value = "Itga7"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 30.09047233333333
$ws.Range("H17").Value = 90.271417
$ws.Range("I17").Value = 0.148170154805508
$ws.Range("J17").Value = 0.148170154805508
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 47.807061
$ws.Range("N17").Value = 143.421183
$ws.Range("O17").Value = 0.8218802864734033
$ws.Range("P17").Value = 0.8218802864734033
$ws.Range("Q17").Value = 1438.537046358479
$ws.Range("R17").Value = 12946.83341722631
$ws.Range("S17").Value = 0.1217781292783595
$ws.Range("T17").Value = 0.1217781292783595
